$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F20").Value = 238
$ws1.Range("F22").Value = 369
$ws1.Range("F27").Value = 1048
$ws1.Range("F30").Value = 24
$ws1.Range("F38").Value = 1920
$ws1.Range("F39").Value = 4347
$ws1.Range("F47").Value = 22

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F19").Value = 238
$ws4.Range("F23").Value = 370
$ws4.Range("F26").Value = 1048
$ws4.Range("F31").Value = 24
$ws4.Range("F42").Value = 22
